$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32
$ws.Range('A32').Value = 112176088
$ws.Range('AA32').NumberFormat = '@'
$ws.Range('AA32').Value = '2023-06-27'
$ws.Range('AX32').Value = 'Pekka Bader, Anna-Maria Eriksson'
$ws.Range('B32').Value = 89559
$ws.Range('E32').Value = 5442
$ws.Range('F32').Value = 'Tallticka'
$ws.Range('G32').Value = 'Porodaedalea pini'
$ws.Range('H32').Value = '(Brot.) Murrill'
$ws.Range('Q32').Value = 602859
$ws.Range('R32').Value = 7030591
$ws.Range('Y32').NumberFormat = '@'
$ws.Range('Y32').Value = '2023-06-27'

# Row 33
$ws.Range('A33').Value = 112176102
$ws.Range('AA33').NumberFormat = '@'
$ws.Range('AA33').Value = '2023-06-14'
$ws.Range('AX33').Value = 'Pekka Bader'
$ws.Range('B33').Value = 89820
$ws.Range('E33').Value = 658
$ws.Range('F33').Value = 'Rosenticka'
$ws.Range('G33').Value = 'Rhodofomes roseus'
$ws.Range('H33').Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Range('Q33').Value = 602660
$ws.Range('R33').Value = 7030716
$ws.Range('Y33').NumberFormat = '@'
$ws.Range('Y33').Value = '2023-06-14'

# Row 34
$ws.Range('A34').Value = 112176093
$ws.Range('AA34').NumberFormat = '@'
$ws.Range('AA34').Value = '2023-06-22'
$ws.Range('B34').Value = 6202
$ws.Range('D34').Value = 'LC'
$ws.Range('E34').Value = 105336
$ws.Range('F34').Value = 'Vanlig flatbagge'
$ws.Range('G34').Value = 'Peltis ferruginea'
$ws.Range('H34').Value = '(Linnaeus, 1758)'
$ws.Range('Q34').Value = 602865
$ws.Range('R34').Value = 7030578
$ws.Range('Y34').NumberFormat = '@'
$ws.Range('Y34').Value = '2023-06-22'

# Row 35
$ws.Range('A35').Value = 112176069
$ws.Range('AR35').Value = ''
$ws.Range('B35').Value = 18534
$ws.Range('D35').Value = 'EN'
$ws.Range('E35').Value = 101797
$ws.Range('F35').Value = 'Karelsk barkfluga'
$ws.Range('G35').Value = 'Xylomya czekanovskii'
$ws.Range('H35').Value = 'Pleske, 1925'
$ws.Range('K35').Value = 'puppa'
$ws.Range('Q35').Value = 602775
$ws.Range('R35').Value = 7030644

# Row 36
$ws.Range('A36').Value = 112176096
$ws.Range('AA36').NumberFormat = '@'
$ws.Range('AA36').Value = '2023-06-14'
$ws.Range('AC36').Value = 'larv 20-25 mm'
$ws.Range('AR36').Value = ""
$ws.Range('AX36').Value = 'Pekka Bader'
$ws.Range('B36').Value = 12450
$ws.Range('E36').Value = 101692
$ws.Range('F36').Value = 'Större barkplattbagge'
$ws.Range('G36').Value = 'Pytho kolwensis'
$ws.Range('H36').Value = 'Sahlberg, 1833'
$ws.Range('K36').Value = 'larv/nymf'
$ws.Range('Q36').Value = 602869
$ws.Range('R36').Value = 7030590
$ws.Range('Y36').NumberFormat = '@'
$ws.Range('Y36').Value = '2023-06-14'

# Row 37
$ws.Range('A37').Value = 112176087
$ws.Range('AA37').NumberFormat = '@'
$ws.Range('AA37').Value = '2023-06-27'
$ws.Range('AX37').Value = 'Pekka Bader, Anna-Maria Eriksson'
$ws.Range('B37').Value = 89724
$ws.Range('D37').Value = 'VU'
$ws.Range('E37').Value = 48
$ws.Range('F37').Value = 'Lappticka'
$ws.Range('G37').Value = 'Amylocystis lapponica'
$ws.Range('H37').Value = '(Romell) Singer'
$ws.Range('Q37').Value = 602806
$ws.Range('R37').Value = 7030689
$ws.Range('Y37').NumberFormat = '@'
$ws.Range('Y37').Value = '2023-06-27'

# Row 38
$ws.Range('A38').Value = 112176074
$ws.Range('AA38').NumberFormat = '@'
$ws.Range('AA38').Value = '2023-06-27'
$ws.Range('AC38').Value = ""
$ws.Range('AX38').Value = 'Pekka Bader, Anna-Maria Eriksson'
$ws.Range('B38').Value = 98934
$ws.Range('D38').Value = 'LC'
$ws.Range('E38').Value = 1365
$ws.Range('F38').Value = 'Lappranunkel'
$ws.Range('G38').Value = 'Coptidium lapponicum'
$ws.Range('H38').Value = '(L.) Tzvelev'
$ws.Range('K38').Value = ""
$ws.Range('Q38').Value = 602642
$ws.Range('R38').Value = 7030561
$ws.Range('Y38').NumberFormat = '@'
$ws.Range('Y38').Value = '2023-06-27'

# Row 39
$ws.Range('A39').Value = 112176095
$ws.Range('B39').Value = 89979
$ws.Range('D39').Value = 'VU'
$ws.Range('E39').Value = 1209
$ws.Range('F39').Value = 'Rynkskinn'
$ws.Range('G39').Value = 'Phlebia centrifuga'
$ws.Range('H39').Value = 'P.Karst.'
$ws.Range('Q39').Value = 602796
$ws.Range('R39').Value = 7030566

# Row 40
$ws.Range('A40').Value = 112176108
$ws.Range('AA40').NumberFormat = '@'
$ws.Range('AA40').Value = '2023-06-14'
$ws.Range('AX40').Value = 'Pekka Bader'
$ws.Range('B40').Value = 89820
$ws.Range('D40').Value = 'NT'
$ws.Range('E40').Value = 658
$ws.Range('F40').Value = 'Rosenticka'
$ws.Range('G40').Value = 'Rhodofomes roseus'
$ws.Range('H40').Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Range('Q40').Value = 602831
$ws.Range('R40').Value = 7030665
$ws.Range('Y40').NumberFormat = '@'
$ws.Range('Y40').Value = '2023-06-14'
